$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 3 through 17 (the old trailing data rows), keeping header (row 1)
# and the first data row (row 2).
$ws.Range("A3:B17").EntireRow.Delete()

# Update the remaining data row with the new (bugfixed) values.
$ws.Range("A2").Value = 39400
$ws.Range("B2").Value = 5.361718827437545
